# Generate Report for Handoff
#
# The "9afdfc77-1c64-4270-878e-ca47b051a6ea" file's status moves from
# "Handed back: in sync with en-US" to "Ready for handoff" and its latest
# handoff timestamp is refreshed, on the Overview sheet as well as on each
# per-language (zh-cn / de-de) detail sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("D3").Value = "2016-35-17 12:35:36"

# --- zh-cn detail sheet ------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("E3").Value = "2016-03-17 12:35:32"

# --- de-de detail sheet ------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("E3").Value = "2016-03-17 12:35:36"
